# Add 2022-Q3 data:
#  - insert a new worksheet "2022-Q3" right after "总计", populate it with the
#    fund-holding data for the new quarter
#  - update the "总计" (summary) sheet so that a new row for 2022-Q3 is
#    inserted above the existing 2022-Q2 / 2022-Q1 rows

$wb = $excel.ActiveWorkbook

function Set-TextValue($rng, $val) {
    # Force the value to be stored as text (matches source data where
    # numeric-looking figures such as "0.76" or "002504" are kept as strings)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet right after "总计"
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $total)
$q3.Name = "2022-Q3"

# Match the page margins used by the other data sheets (0.75in/1in/0.5in)
$q3.PageSetup.LeftMargin = 54
$q3.PageSetup.RightMargin = 54
$q3.PageSetup.TopMargin = 72
$q3.PageSetup.BottomMargin = 72
$q3.PageSetup.HeaderMargin = 36
$q3.PageSetup.FooterMargin = 36

# Reuse the header/style formatting already used elsewhere in the workbook
$total.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$q3.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
Set-TextValue $q3.Range("B2") "002504"
Set-TextValue $q3.Range("C2") "鹏华金鼎灵活配置混合A"
Set-TextValue $q3.Range("D2") "0.76"
Set-TextValue $q3.Range("E2") "92.38"
Set-TextValue $q3.Range("F2") "4.59"
Set-TextValue $q3.Range("G2") "0.0349"
$q3.Range("H2").Value = 8

$q3.Range("A3").Value = 1
Set-TextValue $q3.Range("B3") "002505"
Set-TextValue $q3.Range("C3") "鹏华金鼎灵活配置混合C"
Set-TextValue $q3.Range("D3") "0.19"
Set-TextValue $q3.Range("E3") "92.38"
Set-TextValue $q3.Range("F3") "4.59"
Set-TextValue $q3.Range("G3") "0.0087"
$q3.Range("H3").Value = 8

# ---------------------------------------------------------------------------
# 2. Update the "总计" sheet: push the existing 2022-Q2 / 2022-Q1 rows down
#    by one row and insert the new 2022-Q3 figures at the top of the table
# ---------------------------------------------------------------------------
$ws = $total

$oldB2 = $ws.Range("B2").Value()
$oldC2 = $ws.Range("C2").Value()
$oldD2 = $ws.Range("D2").Value()
$oldB3 = $ws.Range("B3").Value()
$oldC3 = $ws.Range("C3").Value()
$oldD3 = $ws.Range("D3").Value()

# New row 4 needs the same formatting as the existing index column cells
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# row 3 (old 2022-Q1) -> row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = $oldB3
$ws.Range("C4").Value = $oldC3
$ws.Range("D4").Value = $oldD3

# row 2 (old 2022-Q2) -> row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = $oldB2
$ws.Range("C3").Value = $oldC2
$ws.Range("D3").Value = $oldD2

# new 2022-Q3 figures -> row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2022-Q3"
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 0.04

# ---------------------------------------------------------------------------
# 3. Restore the originally active sheet (2022-Q1) as the selected tab
# ---------------------------------------------------------------------------
$wb.Worksheets.Item("2022-Q1").Activate()
